$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.098.56'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.15%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.836.72'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.40%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.85%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6364'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.16%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.003'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.23%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07605'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.19%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2956'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.17%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.88'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.82%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07765'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.39%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.823.38'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.32%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.012'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.00%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6726'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.54%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.44'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.59%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009815'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +7.69%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.123'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.61%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.095.21'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.17%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.60'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.89%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '227.39'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.83%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.234'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.73%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.002'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '160.70'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.82%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1407'
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.552'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.53%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.04'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.37%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.503'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.42%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.131'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.88%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.070'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.80%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.209'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.75%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05400'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.05%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.867'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.45%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7512'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.40%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.143'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.68%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.667'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.70%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.237.61'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.72%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01799'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.94%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.762'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.53%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.625'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.78%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9044'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.32%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.003'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '102.72'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.97%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.987.65'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.55%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000124'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.94%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.07'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.83%  '

$ws.Range("E47").Value = '  +0.05%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4097'
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.073'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.65%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.789'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.93%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05782'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.55%  '
